$wb = $excel.ActiveWorkbook

# --- Arkusz8: convert column B (B2:B6) from whole-number percentages to
#     true percentage values, formatted as Percentage (0% / 0.00%). ---
$ws8 = $wb.Worksheets.Item(8)

# B2:B4 -> 0% (no decimals)
$ws8.Range("B2").Value = 0.12
$ws8.Range("B3").Value = 0.17
$ws8.Range("B4").Value = 0.28
$ws8.Range("B2:B4").NumberFormat = "0%"

# B5:B6 -> 0.00% (two decimals)
$ws8.Range("B5").Value = 0.294
$ws8.Range("B6").Value = 0.323
$ws8.Range("B5:B6").NumberFormat = "0.00%"

# Arkusz8 becomes the active sheet/tab, with B2 selected.
$ws8.Activate()
$ws8.Range("B2").Select()
